$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for columns B (TB), C (d2S), D (K), E (IP), G (sum)
# Column F (Win) is unchanged.
$data = @{
    2  = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    3  = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    4  = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    5  = @(0.1554434735375247, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1.705647867635037)
    6  = @(0.1554434735375247, 0.05231270169004087, 3.082599426703578, 0.4998867070740569, 3.790242309005201)
    7  = @(0.7287194209349384, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 2.27892381503245)
    8  = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    9  = @(1.505614041169197, 1.65323645889881,  0.1529057820181812, 0.4998867070740569, 3.811642989160245)
    10 = @(0.06328177979961902, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 1.328214039578707)
    11 = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    12 = @(3.182878228561681, 1.65323645889881,  0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    13 = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    14 = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    15 = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    16 = @(1.505614041169197, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 4.371470058157054)
    17 = @(1.505614041169197, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 4.371470058157054)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G - sum
}

$wb.Save()
